# Applies the "Modified Absentees bug and consolidated output bug" edit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Student Summary"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Student Summary")

# Insert the new "Course Code:" / "Max Marks:" rows (11 & 12), copying the
# formatting of the existing header rows above them (row 10 for A/B, and the
# already-styled blank C10 for C).
$ws1.Range("A10:C10").Copy()
$ws1.Range("A11:C11").PasteSpecial(-4122)
$ws1.Range("A10:C10").Copy()
$ws1.Range("A12:C12").PasteSpecial(-4122)

$ws1.Range("B11").Value = "Course Code:"
$ws1.Range("C11").Value = "DSPE605"
$ws1.Range("B12").Value = "Max Marks:"
$ws1.Range("C12").Value = 40

# Re-word the statistics labels.
$ws1.Range("A17").Value = "Average Marks"
$ws1.Range("A18").Value = "Less Than 40%"
$ws1.Range("A19").Value = "Between 40 % - 75 %"
$ws1.Range("A20").Value = "More than 75%"

# Round the average marks value.
$ws1.Range("B17").Value = 19.04

# ---------------------------------------------------------------------------
# Sheet 2: "Slow Learners" -- remove the "Rajadurai P" row (was row 16).
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Slow Learners")
$ws2.Range("A16:D16").EntireRow.Delete()

# ---------------------------------------------------------------------------
# Sheet 3: "Fast Learners" -- append six newly-identified fast learners.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Fast Learners")

$ws3.Range("A7").Value = 30
$ws3.Range("B7").Value = 2136110031
$ws3.Range("C7").Value = "Hitesh Kumar K A"
$ws3.Range("D7").Value = 28

$ws3.Range("A8").Value = 13
$ws3.Range("B8").Value = 2136110013
$ws3.Range("C8").Value = "Naveena A"
$ws3.Range("D8").Value = 27

$ws3.Range("A9").Value = 31
$ws3.Range("B9").Value = 2136110032
$ws3.Range("C9").Value = "Jaikrishnan V"
$ws3.Range("D9").Value = 27

$ws3.Range("A10").Value = 8
$ws3.Range("B10").Value = 2136110008
$ws3.Range("C10").Value = "Jananika B"
$ws3.Range("D10").Value = 26

$ws3.Range("A11").Value = 9
$ws3.Range("B11").Value = 2136110009
$ws3.Range("C11").Value = "Kalaivani S"
$ws3.Range("D11").Value = 26

$ws3.Range("A12").Value = 25
$ws3.Range("B12").Value = 2136110026
$ws3.Range("C12").Value = "Ajay S"
$ws3.Range("D12").Value = 26
